$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename FluidTx7_* column header values to FluidTx_* (constant renamed for appName)
$ws.Range("B1").Value = "FluidTx_MenuCount"
$ws.Range("D1").Value = "FluidTx_Menu1"
$ws.Range("E1").Value = "FluidTx_Menu2"
$ws.Range("F1").Value = "FluidTx_Menu3"

# Clear the lingering cell selection stored in the sheet view
$ws.Range("A1").Select()
